$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Writing a date-shaped string (dd-mm-yyyy) straight into a General
    # cell makes Excel auto-convert it to a date serial. Route it through
    # a text formula + paste-values so it lands as plain text without
    # picking up any new cell style.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 194 ("07-10-2021") - fill in the previously-empty numeric columns
$ws.Cells.Item(194, 2).Value = 50000
$ws.Cells.Item(194, 3).Value = 1.5
$ws.Cells.Item(194, 4).Value = 1.5
$ws.Cells.Item(194, 5).Value = 1.5
$ws.Cells.Item(194, 6).Value = 3
$ws.Cells.Item(194, 7).Value = 1.5

# Row 195 - new entry for "08-10-2021"
Set-TextValue $ws.Cells.Item(195, 1) "08-10-2021"
$ws.Cells.Item(195, 2).Value = 60000
$ws.Cells.Item(195, 3).Value = 1.5
$ws.Cells.Item(195, 4).Value = 1.5
$ws.Cells.Item(195, 5).Value = 1.5
$ws.Cells.Item(195, 6).Value = 6
$ws.Cells.Item(195, 7).Value = 1.5

# Row 196 - new entry for "12-10-2021"
Set-TextValue $ws.Cells.Item(196, 1) "12-10-2021"
$ws.Cells.Item(196, 7).Value = 1.5

$excel.CutCopyMode = $false
